$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.628.01"
$ws.Range("E2").Value = "  -1.79%  "

$ws.Range("D3").Value = "3.416.37"
$ws.Range("E3").Value = "  -2.53%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.28%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.05%  "

$ws.Range("D9").Value = "3.416.82"
$ws.Range("E9").Value = "  -2.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.93%  "

$ws.Range("E11").Value = "  -3.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.436"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.35%  "

$ws.Range("D13").Value = "4.006.69"
$ws.Range("E13").Value = "  -2.49%  "

$ws.Range("E14").Value = "  -0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.88%  "

$ws.Range("D17").Value = "64.666.72"
$ws.Range("E17").Value = "  -1.73%  "

$ws.Range("D18").Value = "3.431.09"
$ws.Range("E18").Value = "  -2.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.546"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.27%  "

$ws.Range("E24").Value = "  +0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "

$ws.Range("E28").Value = "  -1.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.17%  "

$ws.Range("E32").Value = "  -3.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.54%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.72%  "

$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.28%  "

$ws.Range("E37").Value = "  -3.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0750"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.28%  "

$ws.Range("D39").Value = "2.883.51"
$ws.Range("E39").Value = "  -6.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("E43").Value = "  -0.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0315"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.771"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.75%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "316.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.20%  "

$ws.Range("E51").Value = "  -2.40%  "
